$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D to make room for "Status" (shifts D..H -> E..I)
$ws.Columns.Item(4).Insert()

# --- Header row ---
$ws.Cells.Item(1, 1).Value = "ISIN"
$ws.Cells.Item(1, 2).Value = "Stock Name"
$ws.Cells.Item(1, 3).Value = "Mutual Fund"
$ws.Cells.Item(1, 4).Value = "Status"
$ws.Cells.Item(1, 5).Value = "Jan_2026"
$ws.Cells.Item(1, 6).Value = "Dec_2025"
$ws.Cells.Item(1, 7).Value = "Oct_2025"
$ws.Cells.Item(1, 8).Value = "MoM"
$ws.Cells.Item(1, 9).Value = "QoQ"

# --- Data rows ---
# Row 2: Piramal Finance Ltd
$ws.Cells.Item(2, 1).Value = "INE202B01038"
$ws.Cells.Item(2, 2).Value = "Piramal Finance Ltd"
$ws.Cells.Item(2, 3).Value = "quant BFSI Fund"
$ws.Cells.Item(2, 4).Value = "Adding Consistently"
$ws.Cells.Item(2, 5).Value = 9.314091
$ws.Cells.Item(2, 6).Value = 8.685552
$ws.Cells.Item(2, 7).Value = 9.216153
$ws.Cells.Item(2, 8).Value = 0.628539
$ws.Cells.Item(2, 9).Value = 0.09793799999999919

# Row 3: Shriram Finance Limited
$ws.Cells.Item(3, 1).Value = "INE721A01047"
$ws.Cells.Item(3, 2).Value = "Shriram Finance Limited"
$ws.Cells.Item(3, 3).Value = "quant BFSI Fund"
$ws.Cells.Item(3, 4).Value = "Fresh Entry"
$ws.Cells.Item(3, 5).Value = 8.767368
$ws.Cells.Item(3, 6).Value = 0.0
$ws.Cells.Item(3, 7).Value = 0.0
$ws.Cells.Item(3, 8).Value = 8.767368
$ws.Cells.Item(3, 9).Value = 8.767368

# Row 4: Capri Global Capital Limited
$ws.Cells.Item(4, 1).Value = "INE180C01042"
$ws.Cells.Item(4, 2).Value = "Capri Global Capital Limited"
$ws.Cells.Item(4, 3).Value = "quant BFSI Fund"
$ws.Cells.Item(4, 4).Value = "Adding"
$ws.Cells.Item(4, 5).Value = 8.001052
$ws.Cells.Item(4, 6).Value = 7.774914
$ws.Cells.Item(4, 7).Value = 9.748983
$ws.Cells.Item(4, 8).Value = 0.2261379999999997
$ws.Cells.Item(4, 9).Value = -1.747931000000001

# Row 5: HDFC Bank Limited
$ws.Cells.Item(5, 1).Value = "INE040A01034"
$ws.Cells.Item(5, 2).Value = "HDFC Bank Limited"
$ws.Cells.Item(5, 3).Value = "quant BFSI Fund"
$ws.Cells.Item(5, 4).Value = "Adding Consistently"
$ws.Cells.Item(5, 5).Value = 7.336848
$ws.Cells.Item(5, 6).Value = 0.471945
$ws.Cells.Item(5, 7).Value = 0.0
$ws.Cells.Item(5, 8).Value = 6.864903
$ws.Cells.Item(5, 9).Value = 7.336848

# Row 6: LIC Housing Finance Ltd
$ws.Cells.Item(6, 1).Value = "INE115A01026"
$ws.Cells.Item(6, 2).Value = "LIC Housing Finance Ltd"
$ws.Cells.Item(6, 3).Value = "quant BFSI Fund"
$ws.Cells.Item(6, 4).Value = "Reducing Consistently"
$ws.Cells.Item(6, 5).Value = 6.249933
$ws.Cells.Item(6, 6).Value = 6.340376
$ws.Cells.Item(6, 7).Value = 7.576282
$ws.Cells.Item(6, 8).Value = -0.09044299999999961
$ws.Cells.Item(6, 9).Value = -1.326349

# Row 7: Kotak Mahindra Bank Limited
$ws.Cells.Item(7, 1).Value = "INE237A01036"
$ws.Cells.Item(7, 2).Value = "Kotak Mahindra Bank Limited"
$ws.Cells.Item(7, 3).Value = "quant BFSI Fund"
$ws.Cells.Item(7, 4).Value = "Fresh Entry"
$ws.Cells.Item(7, 5).Value = 5.747081
$ws.Cells.Item(7, 6).Value = 0.0
$ws.Cells.Item(7, 7).Value = 0.0
$ws.Cells.Item(7, 8).Value = 5.747081
$ws.Cells.Item(7, 9).Value = 5.747081

# Row 8: Bajaj Finance Limited
$ws.Cells.Item(8, 1).Value = "INE296A01032"
$ws.Cells.Item(8, 2).Value = "Bajaj Finance Limited"
$ws.Cells.Item(8, 3).Value = "quant BFSI Fund"
$ws.Cells.Item(8, 4).Value = "Fresh Entry"
$ws.Cells.Item(8, 5).Value = 5.264682
$ws.Cells.Item(8, 6).Value = 0.0
$ws.Cells.Item(8, 7).Value = 3.67717
$ws.Cells.Item(8, 8).Value = 5.264682
$ws.Cells.Item(8, 9).Value = 1.587512

# Row 9: ICICI Prudential AMC Ltd
$ws.Cells.Item(9, 1).Value = "INE346A01027"
$ws.Cells.Item(9, 2).Value = "ICICI Prudential AMC Ltd"
$ws.Cells.Item(9, 3).Value = "quant BFSI Fund"
$ws.Cells.Item(9, 4).Value = "Adding Consistently"
$ws.Cells.Item(9, 5).Value = 5.021502
$ws.Cells.Item(9, 6).Value = 4.448358
$ws.Cells.Item(9, 7).Value = 0.0
$ws.Cells.Item(9, 8).Value = 0.5731440000000001
$ws.Cells.Item(9, 9).Value = 5.021502

# Row 10: HDFC Asset Management Company Ltd
$ws.Cells.Item(10, 1).Value = "INE127D01025"
$ws.Cells.Item(10, 2).Value = "HDFC Asset Management Company Ltd"
$ws.Cells.Item(10, 3).Value = "quant BFSI Fund"
$ws.Cells.Item(10, 4).Value = "Adding Consistently"
$ws.Cells.Item(10, 5).Value = 4.974659
$ws.Cells.Item(10, 6).Value = 3.728754
$ws.Cells.Item(10, 7).Value = 0.0
$ws.Cells.Item(10, 8).Value = 1.245905
$ws.Cells.Item(10, 9).Value = 4.974659

# Row 11: HDFC Life Insurance Co Ltd
$ws.Cells.Item(11, 1).Value = "INE795G01014"
$ws.Cells.Item(11, 2).Value = "HDFC Life Insurance Co Ltd"
$ws.Cells.Item(11, 3).Value = "quant BFSI Fund"
$ws.Cells.Item(11, 4).Value = "Adding Consistently"
$ws.Cells.Item(11, 5).Value = 4.448754
$ws.Cells.Item(11, 6).Value = 3.723976
$ws.Cells.Item(11, 7).Value = 2.912382
$ws.Cells.Item(11, 8).Value = 0.7247780000000001
$ws.Cells.Item(11, 9).Value = 1.536372

# Row 12: Adani Enterprises Limited
$ws.Cells.Item(12, 1).Value = "INE423A01024"
$ws.Cells.Item(12, 2).Value = "Adani Enterprises Limited"
$ws.Cells.Item(12, 3).Value = "quant BFSI Fund"
$ws.Cells.Item(12, 4).Value = "Reducing"
$ws.Cells.Item(12, 5).Value = 3.418112
$ws.Cells.Item(12, 6).Value = 3.749559
$ws.Cells.Item(12, 7).Value = 0.0
$ws.Cells.Item(12, 8).Value = -0.3314470000000003
$ws.Cells.Item(12, 9).Value = 3.418112

# Row 13: Nippon Life India Asset Management Ltd
$ws.Cells.Item(13, 1).Value = "INE298J01013"
$ws.Cells.Item(13, 2).Value = "Nippon Life India Asset Management Ltd"
$ws.Cells.Item(13, 3).Value = "quant BFSI Fund"
$ws.Cells.Item(13, 4).Value = "Fresh Entry"
$ws.Cells.Item(13, 5).Value = 2.141619
$ws.Cells.Item(13, 6).Value = 0.0
$ws.Cells.Item(13, 7).Value = 0.0
$ws.Cells.Item(13, 8).Value = 2.141619
$ws.Cells.Item(13, 9).Value = 2.141619

# Row 14: ICICI Bank Limited
$ws.Cells.Item(14, 1).Value = "INE090A01021"
$ws.Cells.Item(14, 2).Value = "ICICI Bank Limited"
$ws.Cells.Item(14, 3).Value = "quant BFSI Fund"
$ws.Cells.Item(14, 4).Value = "Fresh Entry"
$ws.Cells.Item(14, 5).Value = 0.489655
$ws.Cells.Item(14, 6).Value = 0.0
$ws.Cells.Item(14, 7).Value = 0.0
$ws.Cells.Item(14, 8).Value = 0.489655
$ws.Cells.Item(14, 9).Value = 0.489655

# Row 15: ICICI Prudential Life Insurance Co Ltd
$ws.Cells.Item(15, 1).Value = "INE726G01019"
$ws.Cells.Item(15, 2).Value = "ICICI Prudential Life Insurance Co Ltd"
$ws.Cells.Item(15, 3).Value = "quant BFSI Fund"
$ws.Cells.Item(15, 4).Value = "Reducing"
$ws.Cells.Item(15, 5).Value = 0.039149
$ws.Cells.Item(15, 6).Value = 2.836314
$ws.Cells.Item(15, 7).Value = 0.0
$ws.Cells.Item(15, 8).Value = -2.797165
$ws.Cells.Item(15, 9).Value = 0.039149

# Row 16: SBI Cards & Payment Services Ltd
$ws.Cells.Item(16, 1).Value = "INE018E01016"
$ws.Cells.Item(16, 2).Value = "SBI Cards & Payment Services Ltd"
$ws.Cells.Item(16, 3).Value = "quant BFSI Fund"
$ws.Cells.Item(16, 4).Value = "Complete Exit"
$ws.Cells.Item(16, 5).Value = 0.0
$ws.Cells.Item(16, 6).Value = 6.086461
$ws.Cells.Item(16, 7).Value = 5.459195
$ws.Cells.Item(16, 8).Value = -6.086461
$ws.Cells.Item(16, 9).Value = -5.459195

# Row 17: Anand Rathi Wealth Limited
$ws.Cells.Item(17, 1).Value = "INE463V01026"
$ws.Cells.Item(17, 2).Value = "Anand Rathi Wealth Limited"
$ws.Cells.Item(17, 3).Value = "quant BFSI Fund"
$ws.Cells.Item(17, 4).Value = "Complete Exit"
$ws.Cells.Item(17, 5).Value = 0.0
$ws.Cells.Item(17, 6).Value = 7.22842
$ws.Cells.Item(17, 7).Value = 8.099463
$ws.Cells.Item(17, 8).Value = -7.22842
$ws.Cells.Item(17, 9).Value = -8.099463

# Row 18: Kotak Mahindra Bank Limited
$ws.Cells.Item(18, 1).Value = "INE237A01028"
$ws.Cells.Item(18, 2).Value = "Kotak Mahindra Bank Limited"
$ws.Cells.Item(18, 3).Value = "quant BFSI Fund"
$ws.Cells.Item(18, 4).Value = "Complete Exit"
$ws.Cells.Item(18, 5).Value = 0.0
$ws.Cells.Item(18, 6).Value = 3.139423
$ws.Cells.Item(18, 7).Value = 0.0
$ws.Cells.Item(18, 8).Value = -3.139423
$ws.Cells.Item(18, 9).Value = 0.0

# Row 19: Canara HSBC Life Insurance Company Ltd
$ws.Cells.Item(19, 1).Value = "INE01TY01017"
$ws.Cells.Item(19, 2).Value = "Canara HSBC Life Insurance Company Ltd"
$ws.Cells.Item(19, 3).Value = "quant BFSI Fund"
$ws.Cells.Item(19, 4).Value = "Complete Exit"
$ws.Cells.Item(19, 5).Value = 0.0
$ws.Cells.Item(19, 6).Value = 10.298525
$ws.Cells.Item(19, 7).Value = 9.253373
$ws.Cells.Item(19, 8).Value = -10.298525
$ws.Cells.Item(19, 9).Value = -9.253373

# Row 20: SBI Life Insurance Company Limited
$ws.Cells.Item(20, 1).Value = "INE123W01016"
$ws.Cells.Item(20, 2).Value = "SBI Life Insurance Company Limited"
$ws.Cells.Item(20, 3).Value = "quant BFSI Fund"
$ws.Cells.Item(20, 4).Value = "Complete Exit"
$ws.Cells.Item(20, 5).Value = 0.0
$ws.Cells.Item(20, 6).Value = 0.0
$ws.Cells.Item(20, 7).Value = 2.959871
$ws.Cells.Item(20, 8).Value = 0.0
$ws.Cells.Item(20, 9).Value = -2.959871

# Row 21: Life Insurance Corporation Of India
$ws.Cells.Item(21, 1).Value = "INE0J1Y01017"
$ws.Cells.Item(21, 2).Value = "Life Insurance Corporation Of India"
$ws.Cells.Item(21, 3).Value = "quant BFSI Fund"
$ws.Cells.Item(21, 4).Value = "Complete Exit"
$ws.Cells.Item(21, 5).Value = 0.0
$ws.Cells.Item(21, 6).Value = 0.0
$ws.Cells.Item(21, 7).Value = 8.637842
$ws.Cells.Item(21, 8).Value = 0.0
$ws.Cells.Item(21, 9).Value = -8.637842

# Row 22: State Bank of India
$ws.Cells.Item(22, 1).Value = "INE062A01020"
$ws.Cells.Item(22, 2).Value = "State Bank of India"
$ws.Cells.Item(22, 3).Value = "quant BFSI Fund"
$ws.Cells.Item(22, 4).Value = "Complete Exit"
$ws.Cells.Item(22, 5).Value = 0.0
$ws.Cells.Item(22, 6).Value = 0.0
$ws.Cells.Item(22, 7).Value = 9.771674
$ws.Cells.Item(22, 8).Value = 0.0
$ws.Cells.Item(22, 9).Value = -9.771674

# Row 23: Bajaj Finserv Ltd.
$ws.Cells.Item(23, 1).Value = "INE918I01026"
$ws.Cells.Item(23, 2).Value = "Bajaj Finserv Ltd."
$ws.Cells.Item(23, 3).Value = "quant BFSI Fund"
$ws.Cells.Item(23, 4).Value = "Complete Exit"
$ws.Cells.Item(23, 5).Value = 0.0
$ws.Cells.Item(23, 6).Value = 1.978361
$ws.Cells.Item(23, 7).Value = 1.028731
$ws.Cells.Item(23, 8).Value = -1.978361
$ws.Cells.Item(23, 9).Value = -1.028731

